# Append the 2025-03-29 price row to every sheet in the workbook.
# Each worksheet holds a Date/Price history table in columns A:B that runs
# through row 27 (2025-03-28). A new row 28 is appended with the next day's
# date and that day's price — which, except for the USD_CNY exchange-rate
# sheet, simply repeats the prior day's price.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-29"

# Sheet name -> new Price (column B) value for row 28.
$prices = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "42"
    "N-type Wafer"              = "1.21"
    "Cell Topcon 183mm"         = "0.303"
    "Module Topcon 183mm"       = "0.1"
    "Silver Rear_side"          = "5,533"
    "Silver Busbar front-side"  = "8,284"
    "Silver finger front-side"  = "8,334"
    "USD_CNY"                   = "7.2817"
}

$xlPasteValues = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues
$xlUp = [Microsoft.Office.Interop.Excel.XlDirection]::xlUp

function Set-PlainTextValue($targetCell, $helperCell, $text) {
    # Every cell in these tables (dates in column A, prices in column B) is
    # stored as literal text, even the numeric-looking ones. Assigning
    # Range.Value directly would let Excel auto-coerce "2025-03-29" into a
    # date serial or "40" / "5,533" into a number. Routing the literal
    # through a quoted formula (="...") and pasting only the computed value
    # back keeps it as plain text without minting a new number-format style
    # on the cell (which a NumberFormat="@" + Style="Normal" round trip
    # would otherwise leave behind).
    $helperCell.Formula = '="' + $text + '"'
    $helperCell.Copy()
    $targetCell.PasteSpecial($xlPasteValues)
    $helperCell.ClearContents()
    $excel.CutCopyMode = 0
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $price = $prices[$ws.Name]
    if ($null -eq $price) { continue }

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
    $newRow = $lastRow + 1

    $dateCell = $ws.Cells.Item($newRow, 1)
    $priceCell = $ws.Cells.Item($newRow, 2)
    # Scratch cell, far away from the data table, used only to build each
    # text value via a formula before pasting it as a static value.
    $helper = $ws.Cells.Item(1, 10)

    Set-PlainTextValue $dateCell $helper $newDate
    Set-PlainTextValue $priceCell $helper $price
}
